$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "89×58=" "90×75="
Replace-Text "84×41=" "14×71="
Replace-Text "54×55=" "21×68="
Replace-Text "18×77=" "92×78="
Replace-Text "31×21=" "14×98="
Replace-Text "90×42=" "73×86="
Replace-Text "83×59=" "92×33="
Replace-Text "15×49=" "31×52="
Replace-Text "25×90=" "43×70="
Replace-Text "24×46=" "33×45="
Replace-Text "83×25=" "91×43="
Replace-Text "87×31=" "17×11="
Replace-Text "11×54=" "54×61="
Replace-Text "47×95=" "57×85="
Replace-Text "16×99=" "74×61="
Replace-Text "59×35=" "57×49="
Replace-Text "24×93=" "38×86="
Replace-Text "17×35=" "42×42="
Replace-Text "43×67=" "35×24="
Replace-Text "33×58=" "41×36="
Replace-Text "40×22=" "39×67="
Replace-Text "46×20=" "93×76="
Replace-Text "58×98=" "99×97="
Replace-Text "51×80=" "89×72="
Replace-Text "22×51=" "69×83="
